$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.825.47"
$ws.Range("E2").Value = "  -0.78%  "
$ws.Range("D3").Value = "3.010.38"
$ws.Range("E3").Value = "  -3.44%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'556.40"
$ws.Range("E5").Value = "  -0.55%  "
$ws.Range("D6").Value = "'153.07"
$ws.Range("E6").Value = "  -5.03%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'0.564"
$ws.Range("E8").Value = "  -2.85%  "
$ws.Range("D9").Value = "3.015.01"
$ws.Range("E9").Value = "  -3.00%  "
$ws.Range("D10").Value = "'0.113"
$ws.Range("E10").Value = "  -1.52%  "
$ws.Range("E11").Value = "  -4.54%  "
$ws.Range("E12").Value = "  -3.14%  "
$ws.Range("D13").Value = "3.530.94"
$ws.Range("E13").Value = "  -3.44%  "
$ws.Range("E14").Value = "  -3.16%  "
$ws.Range("D15").Value = "62.889.57"
$ws.Range("E15").Value = "  -0.68%  "
$ws.Range("D16").Value = "'23.92"
$ws.Range("E16").Value = "  -2.79%  "
$ws.Range("D17").Value = "3.013.16"
$ws.Range("E17").Value = "  -2.98%  "
$ws.Range("E18").Value = "  -1.55%  "
$ws.Range("D19").Value = "'396.59"
$ws.Range("E19").Value = "  +0.23%  "
$ws.Range("D20").Value = "'5.10"
$ws.Range("E20").Value = "  -1.64%  "
$ws.Range("D21").Value = "'11.88"
$ws.Range("E21").Value = "  -4.14%  "
$ws.Range("E22").Value = "  -5.33%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").Value = "'65.17"
$ws.Range("E24").Value = "  -2.87%  "
$ws.Range("D25").Value = "'0.466"
$ws.Range("E25").Value = "  -1.73%  "
$ws.Range("E26").Value = "  -5.69%  "
$ws.Range("D27").Value = "0.0₃0968"
$ws.Range("E27").Value = "  -3.23%  "
$ws.Range("D28").Value = "'8.64"
$ws.Range("E28").Value = "  -0.28%  "
$ws.Range("D29").Value = "'0.996"
$ws.Range("E29").Value = "  -0.27%  "
$ws.Range("E31").Value = "  -1.02%  "
$ws.Range("D32").Value = "'20.45"
$ws.Range("E32").Value = "  -1.63%  "
$ws.Range("D33").Value = "'160.71"
$ws.Range("E33").Value = "  +5.77%  "
$ws.Range("D34").Value = "'4.69"
$ws.Range("E34").Value = "  -1.77%  "
$ws.Range("D35").Value = "'6.04"
$ws.Range("E35").Value = "  -2.73%  "
$ws.Range("E36").Value = "  +0.22%  "
$ws.Range("E37").Value = "  -1.66%  "
$ws.Range("E38").Value = "  -3.50%  "
$ws.Range("D39").Value = "2.489.52"
$ws.Range("E39").Value = "  -8.90%  "
$ws.Range("D40").Value = "'37.63"
$ws.Range("E40").Value = "  -1.46%  "
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").Value = "'22.55"
$ws.Range("E41").Value = "  -2.70%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").Value = "'3.91"
$ws.Range("E42").Value = "  -3.20%  "
$ws.Range("E43").Value = "  -3.86%  "
$ws.Range("E44").Value = "  -2.42%  "
$ws.Range("D45").Value = "'0.999"
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("E46").Value = "  -3.35%  "
$ws.Range("D47").Value = "'5.03"
$ws.Range("E47").Value = "  -7.44%  "
$ws.Range("D48").Value = "'19.96"
$ws.Range("E48").Value = "  -3.45%  "
$ws.Range("E49").Value = "  -2.69%  "
$ws.Range("E50").Value = "  +0.41%  "
$ws.Range("D51").Value = "'263.91"
$ws.Range("E51").Value = "  -5.12%  "
